# Ajout draft mapping f595a2bd5e53be80aa00972cfd76eee4a5f7087b
#
# 1. Metadata sheet: bump the "Date" value to the new publication timestamp.
# 2. Elements sheet: add a new mapping column (AL) for the business
#    specification mapping ("Mapping: Spécification métier vers l'extension
#    ROR HealthcareServicePsychiatricSector"), leaving it blank for every
#    row except Extension.value[x], which is mapped to "utilisation".

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# New column header (AL1), cloning the style of the preceding mapping
# column header (AK1) so the look & feel (bold header style) matches.
$elements.Range("AK1").Copy($elements.Range("AL1"))
$elements.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR HealthcareServicePsychiatricSector"

# Blank data cells for rows 2-5 (Extension, Extension.id, Extension.extension,
# Extension.url), cloned from an existing blank/empty-string styled cell so
# the formatting (vertical align top / wrap text) matches the rest of the
# table.
$elements.Range("D2").Copy($elements.Range("AL2"))
$elements.Range("D3").Copy($elements.Range("AL3"))
$elements.Range("D4").Copy($elements.Range("AL4"))
$elements.Range("D5").Copy($elements.Range("AL5"))

# Row 6 (Extension.value[x]) gets the actual mapping value.
$elements.Range("D6").Copy($elements.Range("AL6"))
$elements.Range("AL6").Value = "utilisation"

# Give the new column a wide, best-fit-like width similar to the other long
# text columns (L, M, N, AC, ...).
$elements.Columns.Item(38).ColumnWidth = 91.1

Write-Output "Applied draft mapping column + date bump"
